$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.310.01"
$ws.Range("E2").Value = "  +8.66%  "
$ws.Range("D3").Value = "1.597.67"
$ws.Range("E3").Value = "  +8.02%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9981"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.71"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3686"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3388"
$ws.Range("E8").Value = "  +10.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.62"
$ws.Range("E9").Value = "  +7.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.138"
$ws.Range("E10").Value = "  +7.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07034"
$ws.Range("E11").Value = "  +5.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.64"
$ws.Range("E13").Value = "  +8.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.915"
$ws.Range("E14").Value = "  +7.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.639"
$ws.Range("E15").Value = "  +7.07%  "
$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9974"
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "1.597.97"
$ws.Range("E17").Value = "  +7.89%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001078"
$ws.Range("E18").Value = "  +4.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06616"
$ws.Range("E19").Value = "  +11.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.97"
$ws.Range("E20").Value = "  +11.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.11"
$ws.Range("E21").Value = "  +10.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.015"
$ws.Range("E22").Value = "  +9.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.74"
$ws.Range("E23").Value = "  +6.30%  "
$ws.Range("D24").Value = "22.321.87"
$ws.Range("E24").Value = "  +8.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.398"
$ws.Range("E25").Value = "  +6.68%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.486"
$ws.Range("E26").Value = "  +16.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.33"
$ws.Range("E27").Value = "  +6.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.54"
$ws.Range("E28").Value = "  +13.02%  "
$ws.Range("D29").Value = "1.777.10"
$ws.Range("E29").Value = "  +8.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.35"
$ws.Range("E30").Value = "  +5.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.158"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.009"
$ws.Range("E32").Value = "  +20.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9442"
$ws.Range("E33").Value = "  +16.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08239"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.614"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.300"
$ws.Range("E36").Value = "  +12.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.612"
$ws.Range("E37").Value = "  +10.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.73"
$ws.Range("E38").Value = "  +12.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06141"
$ws.Range("E39").Value = "  +5.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.237"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02211"
$ws.Range("E41").Value = "  +7.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2023"
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9966"
$ws.Range("E43").Value = "  +2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5900"
$ws.Range("E44").Value = "  +11.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.13"
$ws.Range("E45").Value = "  +7.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.668"
$ws.Range("E46").Value = "  +3.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5686"
$ws.Range("E47").Value = "  +9.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.55"
$ws.Range("E48").Value = "  +5.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.962"
$ws.Range("E49").Value = "  +8.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06809"
$ws.Range("E50").Value = "  +5.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.46"
$ws.Range("E51").Value = "  +8.44%  "
